$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the first empty row right after the existing data (row 90 -> 91)
$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value = "Kindergarden"
$ws.Cells.Item($newRow, 2).Value = "Kindergarden Den Haag Binckhorstlaan"
$ws.Cells.Item($newRow, 3).Value = "BSO"
$ws.Cells.Item($newRow, 4).NumberFormat = "@"
$ws.Cells.Item($newRow, 4).Value = "2024-02-22"
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
